$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.978.43"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "2.215.98"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "255.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.28"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -3.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.84"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0910"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.93"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.44"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.50%  "
$ws.Range("D16").Value = "2.213.07"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.782"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.35%  "
$ws.Range("D18").Value = "42.880.02"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000103"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.33"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.20"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.16"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.24"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -8.43%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.70"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.17"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.36"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.19"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "174.09"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.28"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0850"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +6.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.22"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("E35").Value = "  -1.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0357"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.63%  "
$ws.Range("E37").Value = "  -3.38%  "
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.63"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.10"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.77"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +15.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.30"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.73%  "
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.23"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.04"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0978"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.30"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.452"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.12"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").Value = "2.436.50"
$ws.Range("E51").Value = "  -0.98%  "
